$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Extend the R1 traceability block (PacManController / loadGame) down one
#    row so it also covers the new getSavedGame(ActionEvent) method, then
#    re-merge G19:G24 / H19:H24 (was G19:G23 / H19:H23).
# ---------------------------------------------------------------------------
$ws.Range("G19:G23").UnMerge()
$ws.Range("H19:H23").UnMerge()

$ws.Range("G24").Value = ""
$ws.Range("H24").Value = ""
$ws.Range("I24").Value = "getSavedGame(ActionEvent)"

$ws.Range("G19:G24").Merge()
$ws.Range("H19:H24").Merge()

# The simulator recomputes borders for every row spanned by a fresh merge;
# restore the original "full box" look on the interior rows that must stay
# visually unchanged.
$restoreBox = $ws.Range("G20:G23,H20:H23")
$restoreBox.Borders.Item(7).LineStyle = 1
$restoreBox.Borders.Item(10).LineStyle = 1
$restoreBox.Borders.Item(8).LineStyle = 1
$restoreBox.Borders.Item(9).LineStyle = 1
$restoreBox.HorizontalAlignment = -4108
$restoreBox.VerticalAlignment = -4108
$restoreBox.WrapText = $true

# I24 gets a lighter "sides only" border (matches the new border used for the
# inserted getSavedGame row).
$ws.Range("I24").Borders.Item(7).LineStyle = 1
$ws.Range("I24").Borders.Item(10).LineStyle = 1
$ws.Range("I24").Borders.Item(8).LineStyle = 0
$ws.Range("I24").Borders.Item(9).LineStyle = 0
$ws.Range("I24").HorizontalAlignment = -4108
$ws.Range("I24").VerticalAlignment = -4108
$ws.Range("I24").WrapText = $true

# ---------------------------------------------------------------------------
# 2. Rebuild the R2 block: loadGame(String) keeps its own 2-row group
#    (R2 / PacManController merged over rows 25:26), saveGame moves down.
# ---------------------------------------------------------------------------
$ws.Range("G25").Value = "R2"
$ws.Range("H25").Value = "PacManController"
$ws.Range("I25").Value = "loadGame(String)"

$ws.Range("G26").Value = ""
$ws.Range("H26").Value = ""
$ws.Range("I26").Value = "saveGame(ActionEvent)"

$ws.Range("G25:G26").Merge()
$ws.Range("H25:H26").Merge()

# ---------------------------------------------------------------------------
# 3. R3 (saveGame) shifts down to row 27.
# ---------------------------------------------------------------------------
$ws.Range("G27").Value = "R3"
$ws.Range("H27").Value = "PacManController"
$ws.Range("I27").Value = "saveGame(ActionEvent)"

# ---------------------------------------------------------------------------
# 4. R4 / R5 shift down one row (now blank Clase/Metodo as before).
# ---------------------------------------------------------------------------
$ws.Range("G28").Value = "R4"
$ws.Range("H28").Value = ""
$ws.Range("I28").Value = ""

$ws.Range("G29").Value = "R5"
$ws.Range("H29").Value = ""
$ws.Range("I29").Value = ""

# ---------------------------------------------------------------------------
# 5. R6 (Visualizar Puntajes) now gets its traceability row filled in.
# ---------------------------------------------------------------------------
$ws.Range("G30").Value = "R6"
$ws.Range("H30").Value = "PacManController"
$ws.Range("I30").Value = "showScores(ActionEvent)"

# ---------------------------------------------------------------------------
# 6. R7 (Mover PacMan) grows to a 2-row group: PacMan/movePacMan() on row 31
#    and a new PacManController/updateGame() row on row 32.
# ---------------------------------------------------------------------------
$ws.Range("G31").Value = "R7"
$ws.Range("H31").Value = "PacMan"
$ws.Range("I31").Value = "movePacMan()"

$ws.Range("G32").Value = ""
$ws.Range("H32").Value = "PacManController"
$ws.Range("I32").Value = "updateGame()"

$ws.Range("G31:G32").Merge()

# ---------------------------------------------------------------------------
# 7. Sheet view: scroll position / selection moved.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$sel = $ws.Range("K29")
$sel.Select()
